# Applies price/volume and symbol-order updates from the "Updated symbol list"
# GitHub Actions commit of Tue Dec 13 15:51:43 UTC 2022.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    # Force text storage so numeric-looking strings (e.g. "0.01379")
    # keep their exact characters/trailing zeros instead of becoming numbers.
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

Set-TextValue "D2" "275.21"
Set-TextValue "D3" "22.97"
Set-TextValue "D4" "6.321"
Set-TextValue "D5" "0.06224"
Set-TextValue "D6" "3.647"
Set-TextValue "D7" "6.630"
Set-TextValue "D8" "1.392"
Set-TextValue "D9" "0.8334"
Set-TextValue "D10" "0.01379"
Set-TextValue "D11" "0.1595"
Set-TextValue "D12" "0.08402"
Set-TextValue "D13" "0.03515"
Set-TextValue "D14" "0.03219"
Set-TextValue "D15" "4.085"
Set-TextValue "D16" "0.09283"
Set-TextValue "D17" "0.001641"
Set-TextValue "D19" "0.006355"
Set-TextValue "D20" "0.005708"
Set-TextValue "D23" "3.722"
Set-TextValue "D25" "0.3353"
Set-TextValue "D26" "0.1260"
Set-TextValue "D28" "0.0002706"
Set-TextValue "D41" "0.007115"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1169"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.003344"
$ws.Range("E43").Value = "42CEJICEJI"
Set-TextValue "D45" "0.00006022"
Set-TextValue "D46" "0.0009909"
Set-TextValue "D48" "0.7826"
Set-TextValue "D50" "0.00002402"
Set-TextValue "D51" "0.01241"
